# Update column F (dSF) values for specific rows to match repulled data.
# Mapping of row number -> new value for column F.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -3
    3  = -4
    9  = -10
    11 = -5
    14 = -2
    21 = -1
    22 = 3
    23 = -6
    28 = -2
    29 = -3
    31 = -4
    39 = -3
    42 = -4
    43 = 0
    47 = -7
    48 = -4
    49 = -2
    52 = 3
    55 = 2
    57 = -6
    58 = -1
    59 = -7
    60 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
